$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165; everything from old row 165 downward
# (through old row 203) shifts down to 166-204, matching the diff.
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new weekly record.
$ws.Cells.Item(165, 1).Value = 4
$ws.Cells.Item(165, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(165, 3).Value = "Los Lagos"
$ws.Cells.Item(165, 4).Value = 44543
$ws.Cells.Item(165, 5).Value = 10
$ws.Cells.Item(165, 6).Value = 100112037
$ws.Cells.Item(165, 7).Value = "Cebollín"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 60
$ws.Cells.Item(165, 11).Value = 6000
$ws.Cells.Item(165, 12).Value = 6000
$ws.Cells.Item(165, 13).Value = 6000
$ws.Cells.Item(165, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(165, 15).Value = "Región Metropolitana"
$ws.Cells.Item(165, 16).Value = 167
$ws.Cells.Item(165, 17).Value = 36
$ws.Cells.Item(165, 18).Value = "Hortaliza"
